$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$donorIds = @(
  "0ce5dd49",
  "2c1001cb",
  "37cc37bf",
  "43faa0b9",
  "4abe3e88",
  "50164f59",
  "5cf70f79",
  "5da96769",
  "6ca3e2f6",
  "790a4fcb",
  "802cc63a",
  "85c3ea4d",
  "942dfafb",
  "9bc6ba8c",
  "a2d65af2",
  "a46f1771",
  "ad58f9da",
  "c7d9a301",
  "ce8732ff",
  "d6f1d567",
  "da9326c9",
  "e09ca7bf",
  "ef53a641"
)
$categories = @(
  "Varies too much to say",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "Varies too much to say",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "Varies too much to say",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "A short paragraph (21 – 60 words)",
  "Varies too much to say",
  "One short sentence (≤ 20 words)",
  "Varies too much to say",
  "One short sentence (≤ 20 words)",
  "Varies too much to say",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "One short sentence (≤ 20 words)",
  "Varies too much to say"
)

$timestamp = 45854.65062808384
$question = "q05_prompt_length"

for ($i = 0; $i -lt $donorIds.Count; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $donorIds[$i]
  $ws.Cells.Item($r, 2).Value = $categories[$i]
  $ws.Cells.Item($r, 3).Value = $question
  $ws.Cells.Item($r, 4).Value = $timestamp
  $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Output "done"
